$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.214.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.26%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "4.032.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.59%  "

$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "518.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.47%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.661"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.17%  "

$ws.Range("E8").Value = "  +0.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.756"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.41%  "

$ws.Range("E10").Value = "  +1.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000331"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.54%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +12.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.681.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.65%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.059.16"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.91%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +8.36%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.32%  "

$ws.Range("E18").Value = "  -0.68%  "

$ws.Range("E19").Value = "  -1.81%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.128.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.30%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "437.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.37%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "101.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +16.03%  "

$ws.Range("E23").Value = "  +6.94%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.41%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.60%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.46%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.92%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.82%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.78%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "688.67"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.66%  "

$ws.Range("E32").Value = "  +3.13%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +17.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "68.12"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.443"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.96%  "

$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "41.59"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.33%  "

$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0878"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.59"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +21.54%  "

$ws.Range("E39").Value = "  +2.66%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.11%  "

$ws.Range("E41").Value = "  -0.06%  "

$ws.Range("E42").Value = "  +2.70%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.22"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.80%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.77%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.152"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.46%  "

$ws.Range("E46").Value = "  +5.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.59%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000271"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +21.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.94%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0340"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.31%  "
